# Apply re-pulled "dSF" (column F) values per commit message:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -1
    4  = 4
    5  = 5
    6  = 12
    7  = 2
    8  = 8
    9  = 6
    10 = 1
    11 = -5
    12 = 6
    13 = -3
    14 = -4
    16 = 2
    17 = 7
    18 = -3
    19 = 3
    20 = -1
    22 = -1
    23 = -1
    24 = 2
    25 = -1
    26 = -4
    27 = -1
    28 = -1
    29 = 4
    30 = 4
    31 = 4
    32 = -1
    33 = 4
    34 = -3
    35 = 2
    36 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
